# Refresh cached market-price / profit figures across the Leviathan_Profits
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Values below mirror the
# updated currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for
# each affected leve row, as pulled by the scheduled pricing runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 849.5454999999999
$ws.Range("I2").Value = 996.875
$ws.Range("J2").Value = 456.66666
$ws.Range("K2").Value = 996.875
$ws.Range("L2").Value = 456.66666
$ws.Range("M2").Value = -883.875
$ws.Range("N2").Value = -682.66666

$ws.Range("H17").Value = 455570.6
$ws.Range("J17").Value = 496877
$ws.Range("L17").Value = 1490631
$ws.Range("N17").Value = -1490967

$ws.Range("H18").Value = 4309
$ws.Range("I18").Value = 4309
$ws.Range("K18").Value = 4309
$ws.Range("M18").Value = -4025

$ws.Range("H19").Value = 1300.32
$ws.Range("I19").Value = 523.6667
$ws.Range("J19").Value = 3297.4285
$ws.Range("K19").Value = 523.6667
$ws.Range("L19").Value = 3297.4285
$ws.Range("M19").Value = -348.6667
$ws.Range("N19").Value = -3647.4285

$ws.Range("H40").Value = 4999.75
$ws.Range("J40").Value = 5999.6665
$ws.Range("L40").Value = 5999.6665
$ws.Range("N40").Value = -6349.6665

$ws.Range("H51").Value = 11908834
$ws.Range("I51").Value = 4499.6665
$ws.Range("J51").Value = 20837084
$ws.Range("K51").Value = 4499.6665
$ws.Range("L51").Value = 20837084
$ws.Range("M51").Value = -4015.6665
$ws.Range("N51").Value = -20838052

$ws.Range("H62").Value = 46324.418
$ws.Range("I62").Value = 87518.25
$ws.Range("K62").Value = 87518.25
$ws.Range("M62").Value = -86894.25

$ws.Range("H65").Value = 46324.418
$ws.Range("I65").Value = 87518.25
$ws.Range("K65").Value = 437591.25
$ws.Range("M65").Value = -434471.25

$ws.Range("H74").Value = 3665.5386
$ws.Range("I74").Value = 2236
$ws.Range("K74").Value = 2236
$ws.Range("M74").Value = -1300

$ws.Range("H76").Value = 3926.3333
$ws.Range("I76").Value = 3489.5
$ws.Range("J76").Value = 4800
$ws.Range("K76").Value = 3489.5
$ws.Range("L76").Value = 4800
$ws.Range("M76").Value = -3174.5
$ws.Range("N76").Value = -5430

$ws.Range("H77").Value = 3665.5386
$ws.Range("I77").Value = 2236
$ws.Range("K77").Value = 11180
$ws.Range("M77").Value = -6500

$ws.Range("H79").Value = 3926.3333
$ws.Range("I79").Value = 3489.5
$ws.Range("J79").Value = 4800
$ws.Range("K79").Value = 3489.5
$ws.Range("L79").Value = 4800
$ws.Range("M79").Value = -2397.5
$ws.Range("N79").Value = -6984

$ws.Range("H100").Value = 2704.3125
$ws.Range("I100").Value = 1927
$ws.Range("J100").Value = 3999.8333
$ws.Range("K100").Value = 1927
$ws.Range("L100").Value = 3999.8333
$ws.Range("M100").Value = -1386
$ws.Range("N100").Value = -5081.8333

$ws.Range("H106").Value = 11772.546
$ws.Range("J106").Value = 22604.6
$ws.Range("L106").Value = 22604.6
$ws.Range("N106").Value = -23866.6

$ws.Range("H112").Value = 1856.0834
$ws.Range("J112").Value = 2382.2856
$ws.Range("L112").Value = 7146.8568
$ws.Range("N112").Value = -9362.856800000001

$ws.Range("H133").Value = 89942.75
$ws.Range("J133").Value = 89942.75
$ws.Range("L133").Value = 89942.75
$ws.Range("N133").Value = -100062.75

$ws.Range("H138").Value = 2122.6316
$ws.Range("I138").Value = 1511.3077
$ws.Range("J138").Value = 2635.3547
$ws.Range("K138").Value = 4533.9231
$ws.Range("L138").Value = 7906.0641
$ws.Range("M138").Value = 606.0769
$ws.Range("N138").Value = -18186.0641

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24343.154
$ws.Range("I32").Value = 4977.51
$ws.Range("J32").Value = 165435.72
$ws.Range("K32").Value = 4977.51
$ws.Range("L32").Value = 165435.72
$ws.Range("M32").Value = -4690.51
$ws.Range("N32").Value = -166009.72

$ws.Range("H53").Value = 13613
$ws.Range("I53").Value = 7919.5
$ws.Range("K53").Value = 7919.5
$ws.Range("M53").Value = -7237.5

$ws.Range("H61").Value = 4003.1614
$ws.Range("I61").Value = 4080.8262
$ws.Range("J61").Value = 3779.875
$ws.Range("K61").Value = 4080.8262
$ws.Range("L61").Value = 3779.875
$ws.Range("M61").Value = -3868.8262
$ws.Range("N61").Value = -4203.875

$ws.Range("H74").Value = 1620.6666
$ws.Range("I74").Value = 1405.875
$ws.Range("J74").Value = 2308
$ws.Range("K74").Value = 1405.875
$ws.Range("L74").Value = 2308
$ws.Range("M74").Value = -531.875
$ws.Range("N74").Value = -4056

$ws.Range("H77").Value = 1620.6666
$ws.Range("I77").Value = 1405.875
$ws.Range("J77").Value = 2308
$ws.Range("K77").Value = 7029.375
$ws.Range("L77").Value = 11540
$ws.Range("M77").Value = -2661.375
$ws.Range("N77").Value = -20276

$ws.Range("H102").Value = 2799.6155
$ws.Range("I102").Value = 2397.818
$ws.Range("K102").Value = 2397.818
$ws.Range("M102").Value = -775.8180000000002

$ws.Range("H136").Value = 4003.1614
$ws.Range("I136").Value = 4080.8262
$ws.Range("J136").Value = 3779.875
$ws.Range("K136").Value = 12242.4786
$ws.Range("L136").Value = 11339.625
$ws.Range("M136").Value = -9692.4786
$ws.Range("N136").Value = -16439.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4711.3
$ws.Range("I99").Value = 4711.3
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4711.3
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3213.3
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 1481.1428
$ws.Range("I105").Value = 1059.5
$ws.Range("J105").Value = 4011
$ws.Range("K105").Value = 1059.5
$ws.Range("L105").Value = 4011
$ws.Range("M105").Value = 687.5
$ws.Range("N105").Value = -7505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 283.25
$ws.Range("J22").Value = 129
$ws.Range("L22").Value = 129
$ws.Range("N22").Value = -829

$ws.Range("H105").Value = 1816.4286
$ws.Range("I105").Value = 1810
$ws.Range("J105").Value = 1828
$ws.Range("K105").Value = 1810
$ws.Range("L105").Value = 1828
$ws.Range("M105").Value = -63
$ws.Range("N105").Value = -5322

$ws.Range("H108").Value = 49500
$ws.Range("J108").Value = 49500
$ws.Range("L108").Value = 49500
$ws.Range("N108").Value = -57180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 20001572
$ws.Range("J114").Value = 2050.5
$ws.Range("L114").Value = 6151.5
$ws.Range("N114").Value = -12659.5

$ws.Range("H132").Value = 2355.5293
$ws.Range("J132").Value = 2964.7
$ws.Range("L132").Value = 26682.3
$ws.Range("N132").Value = -31742.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3953.64
$ws.Range("J80").Value = 5091.6665
$ws.Range("L80").Value = 5091.6665
$ws.Range("N80").Value = -7087.6665

$ws.Range("H83").Value = 3953.64
$ws.Range("J83").Value = 5091.6665
$ws.Range("L83").Value = 25458.3325
$ws.Range("N83").Value = -35442.3325

$ws.Range("H97").Value = 39666.176
$ws.Range("I97").Value = 50940.54
$ws.Range("J97").Value = 3024.5
$ws.Range("K97").Value = 50940.54
$ws.Range("L97").Value = 3024.5
$ws.Range("M97").Value = -50444.54
$ws.Range("N97").Value = -4016.5

$ws.Range("H102").Value = 1291.5834
$ws.Range("I102").Value = 1357.0952
$ws.Range("K102").Value = 1357.0952
$ws.Range("M102").Value = 264.9048

$ws.Range("H113").Value = 2939.652
$ws.Range("I113").Value = 2242.923
$ws.Range("K113").Value = 2242.923
$ws.Range("M113").Value = -72.92299999999977

$ws.Range("H132").Value = 3266.8667
$ws.Range("I132").Value = 3207.1035
$ws.Range("K132").Value = 9621.3105
$ws.Range("M132").Value = -7091.3105

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4427.778
$ws.Range("I40").Value = 3910
$ws.Range("K40").Value = 3910
$ws.Range("M40").Value = -3774

$ws.Range("H46").Value = 25217.611
$ws.Range("I46").Value = 43296.7
$ws.Range("J46").Value = 2618.75
$ws.Range("K46").Value = 43296.7
$ws.Range("L46").Value = 2618.75
$ws.Range("M46").Value = -43108.7
$ws.Range("N46").Value = -2994.75

$ws.Range("H68").Value = 2063.2222
$ws.Range("I68").Value = 1544.8334
$ws.Range("K68").Value = 1544.8334
$ws.Range("M68").Value = -795.8334

$ws.Range("H71").Value = 2063.2222
$ws.Range("I71").Value = 1544.8334
$ws.Range("K71").Value = 7724.166999999999
$ws.Range("M71").Value = -3980.166999999999

$ws.Range("H93").Value = 22121.438
$ws.Range("I93").Value = 1374
$ws.Range("J93").Value = 333333
$ws.Range("K93").Value = 1374
$ws.Range("L93").Value = 333333
$ws.Range("M93").Value = -126
$ws.Range("N93").Value = -335829

$ws.Range("H100").Value = 28914.072
$ws.Range("I100").Value = 5272.857
$ws.Range("K100").Value = 5272.857
$ws.Range("M100").Value = -4731.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11895.583
$ws.Range("I62").Value = 10041.667
$ws.Range("J62").Value = 13749.5
$ws.Range("K62").Value = 10041.667
$ws.Range("L62").Value = 13749.5
$ws.Range("M62").Value = -9417.666999999999
$ws.Range("N62").Value = -14997.5

$ws.Range("H65").Value = 11895.583
$ws.Range("I65").Value = 10041.667
$ws.Range("J65").Value = 13749.5
$ws.Range("K65").Value = 50208.335
$ws.Range("L65").Value = 68747.5
$ws.Range("M65").Value = -47088.335
$ws.Range("N65").Value = -74987.5

$ws.Range("H96").Value = 3798.1428
$ws.Range("I96").Value = 4819.8
$ws.Range("J96").Value = 1244
$ws.Range("K96").Value = 4819.8
$ws.Range("L96").Value = 1244
$ws.Range("M96").Value = -3446.8
$ws.Range("N96").Value = -3990

$ws.Range("H107").Value = 21740256
$ws.Range("I107").Value = 1047.875
$ws.Range("K107").Value = 3143.625
$ws.Range("M107").Value = -1223.625
